$d = $word.ActiveDocument

# --- Reorder the intro: the date paragraph moves to the top, the
#     "Opmerkingen ..." heading moves below it (and the blank paragraph
#     that used to separate them is dropped). ---

# 1. Remove the empty paragraph that sits between the heading and the date.
$d.Paragraphs.Item(2).Range.Delete()

# 2. Cut the heading paragraph (style + text + its paragraph mark).
$headingRange = $d.Paragraphs.Item(1).Range
$headingRange.Cut() | Out-Null

# 3. Paste it back right after the date paragraph (which is now first).
$dateRange = $d.Paragraphs.Item(1).Range
$pasteTarget = $d.Range($dateRange.End, $dateRange.End)
$pasteTarget.Paste()

# 4. Re-type the heading text as a single run (this also clears the
#    spell-check proofErr markers that used to bracket "NoPressure").
$d.Content.Find.Execute("Opmerkingen naar aanleiding van testen NoPressure mijn Inzet", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Opmerkingen naar aanleiding van testen NoPressure mijn Inzet", 2)

# --- Add the "Openstaande bugs:" line in front of the bookmark. ---
$bookmarkPara = $d.Paragraphs.Item(3)
$insertPoint = $d.Range($bookmarkPara.Range.Start, $bookmarkPara.Range.Start)
$insertPoint.InsertBefore("Openstaande bugs:")

# --- Append two blank paragraphs at the end of the document body. ---
$d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertParagraphAfter()
$d.Paragraphs.Item($d.Paragraphs.Count).Range.InsertParagraphAfter()
